# "Replacing "meat bag" with "humanoid" for professional reasons"
#
# Sheet1 ("Convos - HelloWorld") and Sheet2 ("Convos - HelloWorldUtt") both
# contain a bot reply cell with the text "Hello, meat bag! How can I help you ?".
# Update every occurrence to say "humanoid" instead.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$newText = "Hello, humanoid! How can I help you ?"

# --- Sheet1!C3 : plain find & replace -------------------------------------
$ws1.Range("C3").Value = $newText

# --- Sheet1!C6 : retyped with "humanoid" highlighted as its own run -------
$ws1.Range("C6").Value = $newText
$run = $ws1.Range("C6").Characters(8, 8)
$run.Font.Name = "Calibri"
$run.Font.Size = 11

# --- Sheet2!C3 : same reply, also retyped with the run-level formatting ---
$ws2.Range("C3").Value = $newText
$run2 = $ws2.Range("C3").Characters(8, 8)
$run2.Font.Name = "Calibri"
$run2.Font.Size = 11

# --- Selection / active-tab bookkeeping ------------------------------------
# Last thing the editor did was leave sheet1's cursor on C6, then flip to
# sheet2 and leave its cursor on C3 - making sheet2 the active tab on save.
$ws1.Range("C6").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("C3").Select() | Out-Null
